$wb = $excel.ActiveWorkbook

# Add a brand new sheet right after "Harpy" and rename it "Lesser Hydra".
$harpy = $wb.Worksheets.Item("Harpy")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $harpy)
$newSheet.Name = "Lesser Hydra"

# Header row, matching the other hit-location tables in this workbook.
$newSheet.Cells.Item(1, 1).Value = "Location"
$newSheet.Cells.Item(1, 2).Value = "D20"
$newSheet.Cells.Item(1, 3).Value = "Armor"
$newSheet.Cells.Item(1, 4).Value = "HP"

$newSheet.Cells.Item(2, 1).Value = "Body"
$newSheet.Cells.Item(2, 2).Value = "01-02"
$newSheet.Cells.Item(2, 3).Value = 6
$newSheet.Cells.Item(2, 4).Value = 7

$newSheet.Cells.Item(3, 1).Value = "Heads"
$newSheet.Cells.Item(3, 2).Value = "03-20"
$newSheet.Cells.Item(3, 3).Value = 6
$newSheet.Cells.Item(3, 4).Value = 4

# "Harpy" ends up being the selected/active sheet after this edit.
$harpy.Activate()
